$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add hours total for week 1
$ws.Range("B2").Value = "5 uur 30 minuten"

# Row 5: new "week 2" header row, formatted like row 2 (bold week header)
$ws.Range("A2:B2").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)
$ws.Range("A5").Value = "week 2"

# Row 6
$ws.Range("A6").Value = 43514
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("C6").Value = "Ervoor zorgen dat er maximaal 2 markers toegevoegt kunnen worden, de markers een andere opmaak geven"
$ws.Range("B6").Value = "40 minuten"

# Row 7
$ws.Range("A7").Value = 43516
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = "10 minuten"
$ws.Range("C7").Value = "Teken een lijn tussen 2 markers"

# Interesting links rows 9-11
$ws.Range("Q9").Value = "https://developers.google.com/maps/documentation/android-sdk/polygon-tutorial"
$ws.Range("Q10").Value = "https://app-privacy-policy-generator.firebaseapp.com/"
$ws.Range("Q11").Value = "https://firebase.google.com/docs/android/setup"

# Row 8
$ws.Range("A8").Value = 43518
$ws.Range("A8").NumberFormat = "d-mmm"
$ws.Range("B8").Value = "4 uur"
$ws.Range("C8").Value = "Aanmaken project, navigatie toevoegen, schermen toevoegen, verder werken proefproject, polylines blijven updaten wanneer er een marker bij komt, de polyline updaten elke keer dat je locatie verandert"

$ws.Range("C9").Select()
